$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.820.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'1.620.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'211.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'22.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'0.0882"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'1.849.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'1.616.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'3.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "'0.552"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "'64.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "'27.808.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "'226.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "'7.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").Value = "'0.0₃0712"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'9.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "'154.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "'15.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'3.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").Value = "'1.408.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "'3.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").Value = "'0.973"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "'0.843"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "'65.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "'5.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("E45").Value = "  -3.22%  "
$ws.Range("D46").Value = "'1.759.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("D47").Value = "'2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").Value = "'89.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("D50").Value = "'0.0992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.49%  "
$ws.Range("E51").Value = "  -0.48%  "
